# [TM-13] Excel Test Cases Export deals with UUID
# Insert a new "TC_UUID" column into the TEST_CASES sheet, right before
# the existing TC_REFERENCE column (column G), shifting the following
# columns one place to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TEST_CASES")

# Insert a new blank column at G; everything from G onward shifts right.
$ws.Columns.Item(7).Insert()

# Give the new column a header in row 1.
$ws.Cells.Item(1, 7).Value = "TC_UUID"

# Move/restore the active selection as recorded after the edit.
$ws.Range("H12").Select()
